# Update progress (%) values on the "Task List" sheet to reflect task completion.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task List")

$ws.Range("F3").Value = 100
$ws.Range("F4").Value = 100
$ws.Range("F5").Value = 70

$ws.Range("F5").Select()
